$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (S:T) for "latitude" / "longitude" before the
# existing S:W block (recent sale date / generation date / release
# date / extra data 1 / extra data 2), shifting those five columns to
# U:Y.
$ws.Columns("S:T").Insert()

# --- Header row (row 1) ---
$ws.Range("S1").Value = "latitude"
$ws.Range("T1").Value = "longitude"

# --- Latitude / longitude values for each data row ---
$ws.Range("S2").Value = 39.755543000000003
$ws.Range("T2").Value = -105.22110000000001

$ws.Range("S3").Value = 39.751511000000001
$ws.Range("T3").Value = -105.225381

$ws.Range("S4").Value = 39.740510999999998
$ws.Range("T4").Value = -105.17121

$ws.Range("S5").Value = 39.734164
$ws.Range("T5").Value = -105.159808

$ws.Range("S6").Value = 39.733597000000003
$ws.Range("T6").Value = -105.162576

$ws.Range("S7").Value = 39.739111000000001
$ws.Range("T7").Value = -104.984951

$ws.Range("S8").Value = 39.731361
$ws.Range("T8").Value = -104.96074299999999

$ws.Range("S9").Value = 39.741906999999998
$ws.Range("T9").Value = -104.975129

$ws.Range("S10").Value = 39.742130000000003
$ws.Range("T10").Value = -104.996673

$ws.Range("S11").Value = 40.014986
$ws.Range("T11").Value = -105.270546

$ws.Range("S12").Value = 40.007199
$ws.Range("T12").Value = -105.26486800000001

$ws.Range("S13").Value = 39.999355999999999
$ws.Range("T13").Value = -105.26242000000001

$ws.Range("S14").Value = 39.74389
$ws.Range("T14").Value = -105.02010900000001

$ws.Range("S15").Value = 39.748420000000003
$ws.Range("T15").Value = -105.007644

# --- Refresh the worksheet's remembered sort range/key so it covers
# the two newly-inserted columns (A2:T13 -> A2:V13), without actually
# re-ordering the already-arranged rows: apply the sort, then undo
# just that one step so the data stays put but the sort bookkeeping
# (sortState) is left updated.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B2:B13"))
$sortObj.SetRange($ws.Range("A2:V13"))
$sortObj.Header = 2
$sortObj.Apply()
$excel.Undo()

# --- Row heights: header row grew (wrapped "latitude"/"longitude"
# headers), and the two campus-building detail rows grew too.
$ws.Rows("1").RowHeight = 24
$ws.Rows("14").RowHeight = 30
$ws.Rows("15").RowHeight = 30

# --- Column widths for the two new columns: plain default width,
# no bestFit/customWidth like their neighbours.
$ws.Columns("S:T").ColumnWidth = 8.83203125

# --- Update the view: scroll right to the new columns and move the
# active selection to where the user left off.
$ws.Range("U9").Select()
